$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (H1, style index 1:
# bold font + thin border + centered) onto the two new header cells so
# I1/J1 render like the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iVals = @(4,3,7,4,1,1,1,1,1,1,1,1,1,4,8,4,8,6,8,6,7,7,8,7,7,6,8,4,4,1,6,7,6,3,1)
$jVals = @(5,4,7,7,5,4,6,8,7,6,6,7,4,7,8,5,9,9,9,8,9,9,9,8,8,9,8,9,8,3,8,8,9,4,2)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
